# Update the StructureDefinition workbook metadata to reflect the
# LinuxForHealth re-branding (was Alvearie/IBM), bump the version, and
# refresh the publish date, per the gh-pages deployment commit.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/citizen-status"

# --- Metadata sheet -------------------------------------------------
$wsMetadata.Range("B2").Value = $newUrl                       # URL
$wsMetadata.Range("B3").Value = "8.0.0"                       # Version
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"   # Date
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"         # Publisher

# --- Elements sheet ---------------------------------------------------
# The Extension.url element's "Fixed Value" cell mirrors the canonical
# StructureDefinition URL, so it needs to stay in sync too.
$wsElements.Range("Q5").Value = $newUrl
